$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# Update the Startdate / Enddate values in row 2 (leading apostrophe keeps
# the existing "store as text" quote-prefix cell formatting intact).
$ws.Range("H2").Value = "'14/04/2021"
$ws.Range("I2").Value = "'20/05/2021"

# Move the active selection as recorded in the saved workbook.
$ws.Range("J6").Select()
